$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.736.40"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.230.96"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.94"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.13"
$ws.Range("E6").Value = "  -5.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -6.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.79"
$ws.Range("E10").Value = "  -8.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0819"
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("E12").Value = "  -6.68%  "
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("D14").Value = "2.570.42"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.838"
$ws.Range("E15").Value = "  -4.39%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.235.49"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.99"
$ws.Range("E17").Value = "  -3.96%  "
$ws.Range("D18").Value = "43.656.61"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.08"
$ws.Range("E19").Value = "  -8.23%  "
$ws.Range("D20").Value = "0.0₃0963"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.30"
$ws.Range("E21").Value = "  -5.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.30"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.49"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.99"
$ws.Range("E24").Value = "  -6.90%  "
$ws.Range("E25").Value = "  -8.26%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.14"
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("E29").Value = "  -5.87%  "
$ws.Range("E30").Value = "  -9.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.06"
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0830"
$ws.Range("E33").Value = "  -5.43%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.23"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.60"
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  -7.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.116"
$ws.Range("E38").Value = "  -3.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.46"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.54"
$ws.Range("E40").Value = "  -8.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.01"
$ws.Range("E41").Value = "  -11.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0307"
$ws.Range("E42").Value = "  -5.96%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "1.714.57"
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "82.84"
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.194"
$ws.Range("E46").Value = "  -6.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.13"
$ws.Range("E47").Value = "  -4.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.24"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "71.40"
$ws.Range("E49").Value = "  -4.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.63"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.10"
$ws.Range("E51").Value = "  -5.91%  "
